$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (bold, bordered, centered) onto the new H1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add new "Save" header column in H1
$ws.Range("H1").Value = "Save"

# Add new data cell H2 = 0 (plain, unstyled like the other row-2 data cells)
$ws.Range("H2").Value = 0
